# Add two new result sheets, "NC7" and "HC3", at the end of the workbook,
# mirroring the layout/formatting of the existing NC*/HC* result sheets.

$wb = $excel.ActiveWorkbook

# NOTE: this COM-interop shell only reliably binds POSITIONAL function
# parameters — named parameters (`-Foo bar`) don't bind, so every call
# below passes arguments positionally, in declaration order.
function Add-ResultSheet {
    param($Name, $RowLabel, $InVehicle, $AtStop, $Extra, $Total, $TemplateSheetName)

    # Insert the new sheet right after the current last sheet.
    $lastIndex = $wb.Worksheets.Count
    $lastSheet = $wb.Worksheets.Item($lastIndex)
    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $ws.Name = $Name

    # Clone the formatting (styles/borders/fonts) from an existing result
    # sheet so the new sheet reuses the same style records instead of
    # minting new ones. A1 is intentionally left untouched/empty (as on
    # the template sheets), so copy B1:E1 and A2:E2 separately rather
    # than the full A1:E2 block.
    $template = $wb.Worksheets.Item($TemplateSheetName)
    $template.Range("B1:E1").Copy()
    $ws.Range("B1").PasteSpecial(-4122) | Out-Null
    $template.Range("A2:E2").Copy()
    $ws.Range("A2").PasteSpecial(-4122) | Out-Null

    # Now fill in the values/labels for this sheet.
    $ws.Range("B1").Value = "In-vehicle"
    $ws.Range("C1").Value = "At-stop"
    $ws.Range("D1").Value = "Extra"
    $ws.Range("E1").Value = "Total"

    $ws.Range("A2").Value = $RowLabel
    $ws.Range("B2").Value = $InVehicle
    $ws.Range("C2").Value = $AtStop
    $ws.Range("D2").Value = $Extra
    $ws.Range("E2").Value = $Total

    $ws.Range("A1").Select() | Out-Null
}

# NC7 — "No control" figures (matches the other No-control sheets).
Add-ResultSheet "NC7" "No control" 2101.086661275402 12498.70440518066 141.4698672425732 14741.26093369863 "NC6"

# HC3 — "Holding control" figures.
Add-ResultSheet "HC3" "Holding control" 2349.181173559619 12392.65591453851 118.642810481891 14860.47989858003 "HC2"

$wb.Worksheets.Item("NC7").Select() | Out-Null
